# Level1 map edit: clear out the special tile markers (keywall, door, sliding
# door, rock, etc.) that used to sit on top of the level grid, leaving only
# walls (W), empty space (_) and the one remaining "L" ladder/lock tile.
#
# This mirrors the commit that reworked drawing/collision code so sprites for
# things that move away get removed - as part of that the level sheet's
# leftover special-tile letters were cleaned up to plain "_" (except the
# single "L" which stayed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("o_level1")

# Replace the old "special" tokens with "_" (blank/floor) ...
$ws.Range("H2").Value  = "_"
$ws.Range("J4").Value  = "_"
$ws.Range("F7").Value  = "_"
$ws.Range("F9").Value  = "_"
$ws.Range("C10").Value = "_"
$ws.Range("C11").Value = "_"
$ws.Range("I11").Value = "_"

# ... except B11, which keeps its "L" marker.
$ws.Range("B11").Value = "L"

# Move the active selection from M11 to J6.
$ws.Range("J6").Select()
